# implement basic response for events
#
# The task-timer log gets two new entries appended (rows 9-11) and the
# "kitchen" task's running total (C4) is bumped by another hour of elapsed
# time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "kitchen" (row 4) picked up another ~1 hour of elapsed time.
$ws.Range("C4").Value = 13.00138888888889

# Row 9: "new test task" ran for 1 full day.
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("B9").Value = "new test task"
$ws.Range("C9").Value = 1

# Row 10: an event fired without a task name attached; duration logged anyway.
$ws.Range("B8").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = 1

# Row 11: "second test task" ran for about 34.8 minutes.
$ws.Range("B8").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Value = "second test task"
$ws.Range("C11").Value = 0.02416666666666666

$excel.CutCopyMode = $false
